$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create H1 header cell "Save" - copy G1's formatting (bold, border, centered)
# then overwrite the value so the new style entry matches the existing one
# used by the other header cells instead of minting a new, near-duplicate style.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" column values for each data row.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("H8").Value = 0
